# Add an Arabic guidance row beneath the header row of the IRI_Data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IRI_Data")

# Push the existing data row (row 2: SectionCode / Lane / IRI value) down to
# make room for the new guidance row.
$ws.Rows.Item(2).Insert()

# Populate the guidance row with Arabic instructional text.
$ws.Range("A2").Value = "أدخل رمز القطاع/المقطع"
$ws.Range("B2").Value = "أدخل رقم المسار (L1, L2, إلخ)"
$ws.Range("C2").Value = "أدخل قيمة IRI (مؤشر الخشونة الدولي) بوحدة m/km"

# Style the guidance row: italic, dark-gray 9pt text on a light-gray fill,
# right aligned, vertically centered, with wrapped text - and a taller row.
$guidance = $ws.Range("A2:C2")
$guidance.Font.Italic = $true
$guidance.Font.Color = 6710886
$guidance.Font.Size = 9
$guidance.Interior.Color = 15790320
$guidance.HorizontalAlignment = -4152
$guidance.VerticalAlignment = -4108
$guidance.WrapText = $true
$ws.Rows.Item(2).RowHeight = 30

$wb.Save()
